$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.686.23'
$ws.Range('E2').Value = '  +2.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.889.62'
$ws.Range('E3').Value = '  +0.60%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.65'
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4955'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2955'
$ws.Range('E8').Value = '  +1.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06805'
$ws.Range('E9').Value = '  +2.86%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.898.94'
$ws.Range('E10').Value = '  +0.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '17.06'
$ws.Range('E11').Value = '  +1.86%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07313'
$ws.Range('E12').Value = '  +2.00%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '90.72'
$ws.Range('E13').Value = '  +5.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.056'
$ws.Range('E14').Value = '  +4.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6721'
$ws.Range('E15').Value = '  +1.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.700.41'
$ws.Range('E16').Value = '  +2.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007951'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.002'
$ws.Range('E18').Value = '  +0.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.24'
$ws.Range('E19').Value = '  +4.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.144.17'
$ws.Range('E20').Value = '  +1.03%  '
$ws.Range('E21').Value = '  +0.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.854'
$ws.Range('E22').Value = '  +2.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '177.02'
$ws.Range('E23').Value = '  +31.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.058'
$ws.Range('E24').Value = '  +8.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.283'
$ws.Range('E25').Value = '  +1.93%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.67'
$ws.Range('E26').Value = '  +3.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.52'
$ws.Range('E27').Value = '  +10.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.928'
$ws.Range('E28').Value = '  +0.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.379'
$ws.Range('E29').Value = '  +0.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.338'
$ws.Range('E30').Value = '  +4.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08908'
$ws.Range('E31').Value = '  +2.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.022'
$ws.Range('E32').Value = '  +2.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05242'
$ws.Range('E33').Value = '  +4.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7418'
$ws.Range('E34').Value = '  +5.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.135'
$ws.Range('E35').Value = '  +3.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.671'
$ws.Range('E36').Value = '  +0.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01881'
$ws.Range('E37').Value = '  +10.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.696'
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.161'
$ws.Range('E39').Value = '  -1.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9366'
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4349'
$ws.Range('E41').Value = '  +3.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '105.33'
$ws.Range('E42').Value = '  +3.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.804'
$ws.Range('E43').Value = '  -2.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.003'
$ws.Range('E44').Value = '  +0.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.665'
$ws.Range('E45').Value = '  +3.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1360'
$ws.Range('E46').Value = '  +8.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05836'
$ws.Range('E47').Value = '  +2.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.538'
$ws.Range('E48').Value = '  +5.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '33.35'
$ws.Range('E49').Value = '  +2.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3879'
$ws.Range('E50').Value = '  +4.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.381'
$ws.Range('E51').Value = '  +3.47%  '
